$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new Q&A content, written in the same order the author typed it
# (matches the shared-string insertion order captured in the target diff)
$ws.Range('B100').Value = 'Khi sử dụng các thông tin như ID, tên lớp (class), hoặc bộ chọn CSS để tìm kiếm yếu tố mục tiêu trên trang web, bạn cũng có thể gặp phải một số vấn đề. Dưới đây là một số vấn đề thường gặp và cách giải quyết:'
$ws.Range('B102').Value = 'Yếu tố không tìm thấy: Có thể xảy ra tình huống yếu tố mục tiêu không tồn tại trên trang hoặc nó bị ẩn đi. Điều này có thể dẫn đến lỗi khi bạn cố gắng truy cập yếu tố đó. Đảm bảo rằng bạn đã kiểm tra lại thông tin ID, tên lớp, hoặc bộ chọn CSS và đảm bảo rằng yếu tố có thực sự hiện diện trên trang.'
$ws.Range('B103').Value = 'Trùng lặp thông tin: Có thể có nhiều yếu tố trên trang có cùng thông tin ID hoặc tên lớp. Điều này có thể gây ra sự nhầm lẫn và khiến việc xác định yếu tố chính xác trở nên khó khăn. Đảm bảo rằng thông tin bạn sử dụng để xác định yếu tố là duy nhất.'
$ws.Range('B104').Value = 'Thay đổi cấu trúc trang: Nếu cấu trúc trang web thay đổi, các yếu tố có thể di chuyển hoặc bị xóa. Điều này có thể làm hỏng mã của bạn nếu bạn đang sử dụng thông tin tĩnh như ID hoặc tên lớp để tìm kiếm yếu tố. Để giảm thiểu vấn đề này, hãy sử dụng bộ chọn CSS linh hoạt hơn, cho phép bạn tìm kiếm dựa trên cấu trúc và quan hệ của yếu tố.'
$ws.Range('B105').Value = 'Tính nhạy cảm với ngôn ngữ: Các thông tin như ID, tên lớp và bộ chọn CSS có thể bị thay đổi khi ngôn ngữ của trang web thay đổi. Ví dụ, nếu trang web được dịch sang một ngôn ngữ khác, các thông tin như tên lớp có thể thay đổi. Điều này có thể gây ra sự cố khi tìm kiếm yếu tố. Hãy xem xét sử dụng cách khác như xử lý yếu tố dựa trên văn bản hoặc cấu trúc thay vì thông tin tĩnh.'
$ws.Range('B106').Value = 'Độ tin cậy và bảo trì: Sử dụng thông tin như ID và tên lớp có thể dễ dàng thay đổi bởi các nhà phát triển hoặc trong quá trình bảo trì. Điều này có thể làm hỏng mã của bạn. Đảm bảo rằng bạn cập nhật mã của mình khi có sự thay đổi về thông tin này.'
$ws.Range('B107').Value = 'Hiệu suất: Sử dụng các thông tin như ID, tên lớp và bộ chọn CSS để tìm kiếm yếu tố có thể làm cho quá trình tìm kiếm trở nên chậm, đặc biệt trên các trang web phức tạp. Cân nhắc sử dụng cách tìm kiếm khác hoặc tối ưu hóa mã để cải thiện hiệu suất.'
$ws.Range('B109').Value = 'Tóm lại, khi sử dụng thông tin như ID, tên lớp và bộ chọn CSS để tìm kiếm yếu tố trên trang web, hãy cân nhắc các vấn đề tiềm năng và thực hiện các biện pháp đối phó tương ứng để đảm bảo tính chính xác và ổn định của mã của bạn.'
$ws.Range('B112').Value = 'Để đối phó với các vấn đề khi sử dụng thông tin như ID, tên lớp và bộ chọn CSS để tìm kiếm yếu tố trên trang web, bạn có thể áp dụng các biện pháp sau:'
$ws.Range('B122').Value = 'Tóm lại, bằng cách sử dụng các biện pháp như sử dụng bộ chọn CSS linh hoạt, sử dụng thuộc tính không thay đổi, sử dụng xPath, xử lý ngoại lệ, dự phòng thông tin, liên tục kiểm tra và cập nhật mã, cũng như tích hợp với các công cụ kiểm thử tự động, bạn có thể giải quyết các vấn đề liên quan đến việc tìm kiếm yếu tố mục tiêu trên trang web một cách hiệu quả.'
$ws.Range('B114').Value = 'Sử dụng bộ chọn CSS linh hoạt: Thay vì dựa vào thông tin tĩnh như ID hoặc tên lớp, sử dụng các bộ chọn CSS linh hoạt để xác định yếu tố dựa trên vị trí và quan hệ của nó trong cấu trúc HTML. Ví dụ: sử dụng parent > child để tìm kiếm yếu tố con trực tiếp của một yếu tố cha cụ thể.'
$ws.Range('B115').Value = 'Sử dụng thuộc tính không thay đổi: Nếu yếu tố có một thuộc tính không thay đổi như văn bản hoặc giá trị của một thuộc tính khác, bạn có thể sử dụng thông tin này để xác định yếu tố. Ví dụ: sử dụng bộ chọn CSS như [data-id="unique-id"] để tìm kiếm yếu tố dựa trên một thuộc tính tùy chỉnh.'
$ws.Range('B116').Value = 'Sử dụng xPath: XPath là một ngôn ngữ truy vấn sử dụng để xác định các yếu tố trên trang web dựa trên cấu trúc HTML. Xpath cung cấp khả năng xác định yếu tố dựa trên nhiều thông tin khác nhau. Sử dụng xPath để tìm kiếm yếu tố một cách linh hoạt và chính xác.'
$ws.Range('B117').Value = 'Xử lý ngoại lệ và kiểm tra tồn tại: Trước khi truy cập yếu tố, hãy kiểm tra xem yếu tố có tồn tại trên trang hay không. Sử dụng cơ chế xử lý lỗi hoặc hàm kiểm tra sự tồn tại để tránh lỗi khi yếu tố không tìm thấy.'
$ws.Range('B118').Value = 'Dự phòng thông tin: Nếu thông tin như ID hoặc tên lớp thay đổi thường xuyên, bạn có thể dự phòng bằng cách xây dựng nhiều cách tìm kiếm khác nhau để xác định yếu tố. Điều này giúp đảm bảo rằng bạn có nhiều lựa chọn trong trường hợp thông tin thay đổi.'
$ws.Range('B119').Value = 'Liên tục kiểm tra và cập nhật mã: Theo dõi trang web để xem xét liệu có sự thay đổi về cấu trúc, thông tin, hoặc thuộc tính của các yếu tố hay không. Nếu có sự thay đổi, cập nhật mã của bạn để đảm bảo tính ổn định và chính xác.'
$ws.Range('B120').Value = 'Tích hợp với kiểm thử tự động: Sử dụng các framework kiểm thử tự động như Selenium WebDriver để tạo các kịch bản kiểm thử. Đây là các công cụ mạnh mẽ hỗ trợ việc tương tác với trình duyệt và tìm kiếm yếu tố dễ dàng hơn.'

# Update the view to match where the author left the selection/scroll position
$win = $excel.ActiveWindow
$win.ScrollRow = 82
$win.ScrollColumn = 1
$ws.Range('D119').Select() | Out-Null
